$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from 45170
# (2023-09-01) to 45174 (2023-09-05) for every data row (rows 2-70).
$ws.Range("C2:C70").Value = 45174
